# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to reflect the regenerated output at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Mapping of row -> new value for column F
$updates = @{
    5  = 13083
    8  = 515
    9  = 480
    11 = 980
    12 = 13760
    13 = 14341
    22 = 1088
    25 = 5392
    26 = 937
    28 = 309
    29 = 17
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
